# "Generate Report for Handback" — refresh the handoff/handback timestamps
# for the file 4ee89f28-2be6-4d6c-8201-b27c46ddb401.md (row 5 on every sheet)
# after a new handback round-trip completed.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for that file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-10-24 07:29:14"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-10-24 07:29:03"
$wsZhCn.Range("K5").Value = "2016-10-24 07:29:44"

# de-de sheet: same two columns.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-10-24 07:29:14"
$wsDeDe.Range("K5").Value = "2016-10-24 07:30:01"
